# "new meal planner algo" - add a new multi-select question ("Any particular
# food restrictions?") to the "About Your Lifestyle" section of the Quiz
# sheet, as question #5, inserted right before the "One Last Thing" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quiz")

# Make room for the 7 new option rows (old row 91 "One Last Thing / Q1"
# and everything below it shifts down to row 98+).
$ws.Rows.Item(91).Resize(7).Insert()

# Section / question-number columns for the whole new block.
$ws.Range("A91:A97").Value = 3
$ws.Range("B91:B97").Value = "About Your Lifestyle"
$ws.Range("C91:C97").Value = 5

# Question text + type only live on the first row of the option group.
$ws.Range("D91").Value = "Any particular food restrictions?"
$ws.Range("E91").Value = "multiple"

# The seven selectable options.
$ws.Range("F91").Value = "No Meat at all"
$ws.Range("F92").Value = "No Red Meat only"
$ws.Range("F93").Value = "No Fish"
$ws.Range("F94").Value = "No Crustaceans"
$ws.Range("F95").Value = "No Sea Food at all"
$ws.Range("F96").Value = "No Milk & Dairy products"
$ws.Range("F97").Value = "No Eggs"

# Leave the selection where the author's session ended up.
$ws.Activate()
$ws.Range("F106").Select()
